# OLX Monitor 2026-02-17 18:21
# Appends a new scrape block (rows 61-67) to the "PODSUMOWANIE" sheet's
# log table (columns A-H), matching the existing repeating-block layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2026-02-17 18:21:25"

# Each entry: profile, title, price, date, days-listed, url, slug, Warn(red F col)
$entries = @(
    @{ Profile = "poqui"; Title = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda"; Price = 2299; Date = "19.01.2026"; Days = 29; Url = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html"; Slug = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"; Warn = $false },
    @{ Profile = "poqui"; Title = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"; Price = 2499; Date = "28.10.2025"; Days = 112; Url = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"; Slug = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"; Warn = $true },
    @{ Profile = "poqui"; Title = "Przytulny pokój blisko Politechniki – ul. Przytulna"; Price = 599; Date = "10.10.2025"; Days = 130; Url = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"; Slug = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"; Warn = $true },
    @{ Profile = "pokojewlublinie"; Title = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"; Price = 58640; Date = "11.08.2025"; Days = 190; Url = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"; Slug = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"; Warn = $true },
    @{ Profile = "pokojewlublinie"; Title = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"; Price = 12640; Date = "19.01.2026"; Days = 29; Url = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"; Slug = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"; Warn = $false },
    @{ Profile = "dawnypatron"; Title = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."; Price = 730; Date = "20.09.2024"; Days = 515; Url = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"; Slug = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"; Warn = $true },
    @{ Profile = "dawnypatron"; Title = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"; Price = 14690; Date = "05.12.2025"; Days = 74; Url = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"; Slug = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"; Warn = $true }
)

$startRow = 61
$row = $startRow
foreach ($e in $entries) {

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)
    $cellC = $ws.Cells.Item($row, 3)
    $cellD = $ws.Cells.Item($row, 4)
    $cellE = $ws.Cells.Item($row, 5)
    $cellF = $ws.Cells.Item($row, 6)
    $cellG = $ws.Cells.Item($row, 7)
    $cellH = $ws.Cells.Item($row, 8)

    # Text columns that could otherwise be misread as dates/numbers by
    # Excel's input parser ("10.10.2025" etc.) are pre-formatted as Text so
    # they land verbatim as strings, then the format is reset back to
    # General (matching the rest of the sheet) once the literal value is in.
    $cellA.NumberFormat = "@"
    $cellA.Value = $timestamp
    $cellA.Style = "Normal"
    $cellA.HorizontalAlignment = -4131

    $cellB.Value = $e.Profile

    $cellC.NumberFormat = "@"
    $cellC.Value = $e.Title
    $cellC.Style = "Normal"
    $cellC.HorizontalAlignment = -4131

    $cellD.Value = $e.Price
    $cellD.HorizontalAlignment = -4108

    $cellE.NumberFormat = "@"
    $cellE.Value = $e.Date
    $cellE.Style = "Normal"
    $cellE.HorizontalAlignment = -4108

    $cellF.Value = $e.Days
    $cellF.HorizontalAlignment = -4108
    if ($e.Warn) {
        # Long-listed ad -> red highlight, matching the rest of the table.
        $cellF.Font.Name = "Calibri"
        $cellF.Font.Size = 10
        $cellF.Font.Color = 7039999
    }

    $cellG.Value = $e.Url
    $cellH.Value = $e.Slug

    $row = $row + 1
}
